# Hortaliza, Macroferia Regional de Talca - Coliflor
# Weekly refresh: two new daily price observations were added to the
# historical series (rows insert, shifting subsequent rows down).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert new row with date 44748 right after existing row 233 ---
# (this becomes the new row 234, pushing the former rows 234.. down by one)
$ws.Rows.Item(234).Insert(-4121)   # xlShiftDown

$ws.Cells.Item(234, 1).Value2  = 5
$ws.Cells.Item(234, 2).Value2  = "Macroferia Regional de Talca"
$ws.Cells.Item(234, 3).Value2  = "Maule"
$ws.Cells.Item(234, 4).Value2  = 44748
$ws.Cells.Item(234, 5).Value2  = 7
$ws.Cells.Item(234, 6).Value2  = 100112008
$ws.Cells.Item(234, 7).Value2  = "Coliflor"
$ws.Cells.Item(234, 8).Value2  = "Sin especificar"
$ws.Cells.Item(234, 9).Value2  = "Primera"
$ws.Cells.Item(234, 10).Value2 = 3000
$ws.Cells.Item(234, 11).Value2 = 900
$ws.Cells.Item(234, 12).Value2 = 900
$ws.Cells.Item(234, 13).Value2 = 900
$ws.Cells.Item(234, 14).Value2 = "`$/unidad"
$ws.Cells.Item(234, 15).Value2 = "Región del Maule"
$ws.Cells.Item(234, 16).Value2 = 900
$ws.Cells.Item(234, 17).Value2 = 1
$ws.Cells.Item(234, 18).Value2 = "Hortaliza"

# --- Insert new row with date 44747 right after the row now holding date 44265 ---
# (that row is now row 256 after the previous insert, so the new row lands at 257,
# pushing the remaining former rows down by one more)
$ws.Rows.Item(257).Insert(-4121)   # xlShiftDown

$ws.Cells.Item(257, 1).Value2  = 5
$ws.Cells.Item(257, 2).Value2  = "Macroferia Regional de Talca"
$ws.Cells.Item(257, 3).Value2  = "Maule"
$ws.Cells.Item(257, 4).Value2  = 44747
$ws.Cells.Item(257, 5).Value2  = 7
$ws.Cells.Item(257, 6).Value2  = 100112008
$ws.Cells.Item(257, 7).Value2  = "Coliflor"
$ws.Cells.Item(257, 8).Value2  = "Sin especificar"
$ws.Cells.Item(257, 9).Value2  = "Primera"
$ws.Cells.Item(257, 10).Value2 = 3000
$ws.Cells.Item(257, 11).Value2 = 1000
$ws.Cells.Item(257, 12).Value2 = 1000
$ws.Cells.Item(257, 13).Value2 = 1000
$ws.Cells.Item(257, 14).Value2 = "`$/unidad"
$ws.Cells.Item(257, 15).Value2 = "Región del Maule"
$ws.Cells.Item(257, 16).Value2 = 1000
$ws.Cells.Item(257, 17).Value2 = 1
$ws.Cells.Item(257, 18).Value2 = "Hortaliza"

Write-Output "Rows inserted. UsedRange rows now:"
Write-Output $ws.UsedRange.Rows.Count
